$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns
$ws.Range("V1").Value = "Numero Propuesta"
$ws.Range("W1").Value = "Resultado"

# Row 2 changes: Monto (H2) and Tasa Inicial (I2) now stored as text
$ws.Range("H2").Value = "5000"
$ws.Range("I2").Value = "1"
$ws.Range("W2").Value = "ok"

# Row 3 changes: Numero (A3), Monto (H3), Tasa Inicial (I3)
$ws.Range("A3").Value = "2240"
$ws.Range("H3").Value = "5000"
$ws.Range("I3").Value = "1"
$ws.Range("W3").Value = "ok"

# Update selection / view
$ws.Range("B7").Select()
